# fix: added schema settings
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("organisations")

# Headers
$ws.Cells.Item(1, 10).Value = "hasSchema"
$ws.Cells.Item(1, 11).Value = "schemaName"

# Data: row -> (hasSchema, schemaName)
$ws.Cells.Item(2, 10).Value = $true
$ws.Cells.Item(2, 11).Value = "BE1"

$ws.Cells.Item(3, 10).Value = $true
$ws.Cells.Item(3, 11).Value = "CZ1"

$ws.Cells.Item(8, 10).Value = $true
$ws.Cells.Item(8, 11).Value = "DE1"

$ws.Cells.Item(9, 10).Value = $true
$ws.Cells.Item(9, 11).Value = "HU2"

$ws.Cells.Item(12, 10).Value = $true
$ws.Cells.Item(12, 11).Value = "IT2"

$ws.Cells.Item(14, 10).Value = $true
$ws.Cells.Item(14, 11).Value = "IT5"

$ws.Cells.Item(17, 11).Value = "NL1"

$ws.Cells.Item(19, 10).Value = $true
$ws.Cells.Item(19, 11).Value = "NL3"

$ws.Cells.Item(23, 11).Value = "SE1"

# Selection / view changes: organisations becomes the active tab, with K18 selected
$ws.Activate()
$ws.Range("K18").Select()
